# Updates cryptos list price/volume cells (and the Huobi/Lido row swap)
# to match the latest GitHub Actions scrape.
#
# Note: several "Price" cells look like plain decimals (e.g. "9.56").
# Assigning those directly to .Value would make Excel infer a numeric
# type, which would change the underlying cell representation from the
# original plain-text inline string. To keep them as text (matching the
# source data, which mixes European-style "thousands dot" numbers with
# plain-looking decimals - all stored as text), we briefly force the
# cell to Text format ("@") before writing the value, then restore the
# cell style to Normal/General so no stray formatting differences are
# introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.964.21'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.623.45'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('E8').Value = '  -2.41%  '
$ws.Range('E9').Value = '  -3.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.37%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = '1.848.26'
$ws.Range('D13').Value = '1.640.69'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.523'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').Value = '25.938.30'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '0.0₃0737'
$ws.Range('E17').Value = '  -3.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.37%  '
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('E21').Value = '  -2.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('E28').Value = '  -2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.30%  '
$ws.Range('E30').Value = '  -1.64%  '
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.60%  '
$ws.Range('E33').Value = '  -5.71%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('D36').Value = '1.126.82'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.846'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.42%  '
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.771'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').Value = '1.759.15'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.42%  '
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '54.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.80%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.16%  '
